$d = $word.ActiveDocument

# Find the paragraph that talks about Android being built mainly for
# touchscreen devices ("... folosit in principal pe dispozitive cu
# touchscreen ..."). We locate it by content rather than by a fixed
# index so the script is resilient to minor document differences.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*folosit in principal pe dispozitive*touchscreen*") {
        $target = $p.Range
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the target paragraph containing 'folosit in principal'"
}

# Within that paragraph, find the exact character offset of the "in" in
# "folosit in principal" and turn it into the correctly-accented
# Romanian word "in" -> "in" becomes the feminine form? No: the fix is
# simply adding the missing diacritic so "in" (no diacritics) becomes
# "in" with the preposition spelled correctly as "in principal" ->
# "in principal" (i with circumflex), i.e. replace the bare "i" with "i".
$paraText = $target.Text
$needle = "folosit in principal"
$idx = $paraText.IndexOf($needle)
if ($idx -lt 0) {
    throw "Could not locate 'folosit in principal' inside the target paragraph"
}

# Offset of the "i" that starts the word "in" (right after "folosit ").
$iOffset = $idx + "folosit ".Length
$start = $target.Start + $iOffset
$end = $start + 1

$charRange = $d.Range($start, $end)
if ($charRange.Text -ne "i") {
    throw "Unexpected character at computed offset: '$($charRange.Text)'"
}

# Replace the plain "i" with the correctly accented "i" (i-circumflex),
# turning "folosit in principal" into "folosit în principal".
$charRange.Text = [char]0x00EE

Write-Output "Replaced 'in' with 'in' (i-circumflex) in the Android/touchscreen paragraph."
